# The "Export" sheet lists accounts by descending balance ("Saldo").
# Account 004862672 (RENATO) actually has a much smaller balance (7723.73,
# not 52723.73), which moves it several places down the ranking - below
# 004237325 (RICARDO, 9836.3) and above 004397124 (MURYLO, 6148.7).
#
# Concretely: row 3 (004862672/RENATO) needs to move down to just above the
# 004397124/MURYLO row, and rows 4-7 (CARLOS, CINTIA, BRUNO, RICARDO) shift
# up to take its place. RENATO's balance also changes to 7723.73.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert a blank row just above row 8 (MURYLO), i.e. row 7 will
# become free for RENATO once the rows above shift up.
$ws.Range("A8:C8").Insert()

# Move RENATO's row (currently row 3) down into the freshly inserted row 8,
# updating only the name/account columns via Cut (kept as text, no
# precision concerns) while writing the corrected balance directly.
$ws.Range("A3:B3").Cut($ws.Range("A8:B8"))
$ws.Cells.Item(8, 3).Value = 7723.73

# Remove the now-vacated original RENATO row; everything below (CARLOS,
# CINTIA, BRUNO, RICARDO) shifts up by one, landing in rows 3-6, with
# RENATO ending up in row 7 right after RICARDO and before MURYLO.
$ws.Range("A3:C3").Delete()
